# Update "想去人数" (wish-to-go count) values in column F across the
# "展览", "演出" and "全部类型" worksheets, per the upstream gh-pages data
# refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 41
$ws1.Range("F6").Value  = 203
$ws1.Range("F14").Value = 2015
$ws1.Range("F16").Value = 15
$ws1.Range("F18").Value = 467
$ws1.Range("F22").Value = 42
$ws1.Range("F23").Value = 1510
$ws1.Range("F24").Value = 3425
$ws1.Range("F27").Value = 58
$ws1.Range("F28").Value = 1111
$ws1.Range("F29").Value = 98
$ws1.Range("F30").Value = 1799
$ws1.Range("F33").Value = 63
$ws1.Range("F35").Value = 404
$ws1.Range("F39").Value = 49

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 11

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 41
$ws4.Range("F6").Value  = 203
$ws4.Range("F14").Value = 2015
$ws4.Range("F16").Value = 11
$ws4.Range("F17").Value = 15
$ws4.Range("F19").Value = 467
$ws4.Range("F23").Value = 42
$ws4.Range("F24").Value = 1510
$ws4.Range("F25").Value = 3425
$ws4.Range("F28").Value = 58
$ws4.Range("F29").Value = 1111
$ws4.Range("F30").Value = 98
$ws4.Range("F31").Value = 1799
$ws4.Range("F34").Value = 63
$ws4.Range("F36").Value = 404
$ws4.Range("F40").Value = 49
